$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old A1 header cell entirely (it disappears in the target,
#    the header row now starts at B1).
# ------------------------------------------------------------------
$ws.Range("A1").Clear()

# ------------------------------------------------------------------
# 2. Header row: B1:E1 already hold 1..4 (style already applied), add
#    the new F1 = 5 header cell with the same header style.
# ------------------------------------------------------------------
$ws.Range("F1").Value = 5
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Column A (row index helper) now uses the bold/bordered header
#    style on every data row, and holds a simple 1-based counter.
# ------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("B1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4. Column B becomes a plain (unstyled) 0-based counter.
# ------------------------------------------------------------------
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
# B2 previously carried the (now stale) shared-string-cell formatting
# inherited from the old layout - reset formatting to plain/default to
# match target (copy the plain default formatting from D2).
$ws.Range("D2").Copy()
$ws.Range("B2:B5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 5. Columns D and E hold the two numeric amount columns.
# ------------------------------------------------------------------
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 5000
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("E2").Value = 1000
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 25
$ws.Range("E5").Value = 100

# ------------------------------------------------------------------
# 6. Column C (description) and column F (date, kept as literal text
#    rather than an Excel date serial) are both shared-string columns.
#    Fill them row-by-row (description then date, per row) so the
#    shared-string table is built up in the same interleaved order a
#    human filling the sheet in top-to-bottom row order would produce.
#    Force text format on F first so the "mm/dd/yyyy"-shaped string
#    isn't auto-converted into a date value, then strip the number
#    format back off so the cell ends up with the plain/default style.
# ------------------------------------------------------------------
$ws.Range("F2:F5").NumberFormat = "@"

$ws.Range("C2").Value = "Mortage"
$ws.Range("F2").Value = "01/01/2000"

$ws.Range("C3").Value = "Paycheck"
$ws.Range("F3").Value = "01/02/2000"

$ws.Range("C4").Value = "Lunch"
$ws.Range("F4").Value = "01/03/2000"

$ws.Range("C5").Value = "Dog Food"
$ws.Range("F5").Value = "01/04/2000"

$ws.Range("D2").Copy()
$ws.Range("F2:F5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 7. Selection moves to F1 (matches the new header's last cell).
# ------------------------------------------------------------------
[void]$ws.Range("F1").Select()

$excel.CutCopyMode = 0
